# Regenerate save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals.
# This updates the "K" column (column G) values for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 4
    4  = 3
    5  = 0
    6  = 2
    7  = 0
    8  = 2
    9  = 4
    10 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
